$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'sliding pads for exercise'
$ws.Range('A2').Value = 'compression capri men'
$ws.Range('A3').Value = 'basketball pants for boys'
$ws.Range('A4').Value = 'padded baseball pants'
$ws.Range('A5').Value = 'cycling capri pants'
$ws.Range('A6').Value = 'sliding pants baseball mens'
$ws.Range('A7').Value = 'below the knee shorts for men'
$ws.Range('A8').Value = 'mens long basketball shorts below knee'
$ws.Range('A9').Value = 'youth small knee pads basketball'
$ws.Range('A10').Value = 'basketball leg pads'
$ws.Range('A11').Value = 'short baseball pants'
$ws.Range('A12').Value = 'basketball bump'
$ws.Range('A13').Value = 'knee pads running'
$ws.Range('A14').Value = 'football compression shorts with pads'
$ws.Range('A15').Value = 'knee high baseball pants mens'
$ws.Range('A16').Value = 'knee padded compression'
$ws.Range('A17').Value = 'athletic pads'
$ws.Range('A18').Value = 'mens protection pads'
$ws.Range('A19').Value = 'youth knee pads wrestling'
$ws.Range('A20').Value = 'men softball pants'
$ws.Range('A21').Value = 'baseball sliding pants'
$ws.Range('A22').Value = 'boys sliding pants'
$ws.Range('A23').Value = 'youth boys basketball pants'
$ws.Range('A24').Value = 'lacrosse sweat pants'
$ws.Range('A25').Value = 'raceface knee pads'
$ws.Range('A26').Value = 'capri compression pants men'
$ws.Range('A27').Value = 'snowboarding compression pants'
$ws.Range('A28').Value = 'men knee pad pants'
$ws.Range('A29').Value = 'youth basketball tights for boys'
$ws.Range('A30').Value = 'knee pads for basketball youth'
$ws.Range('A31').Value = 'mtn bike knee pads'
$ws.Range('A32').Value = 'baseball padded sliding shorts'
$ws.Range('A33').Value = 'leggings with baseballs'
$ws.Range('A34').Value = 'youth knee pad wrestling'
$ws.Range('A35').Value = '28 basketball'
$ws.Range('A36').Value = 'tights with knee'
$ws.Range('A37').Value = 'compression capris'
$ws.Range('A38').Value = 'padded tights for football'
$ws.Range('A39').Value = 'baseball tights for boys'
$ws.Range('A40').Value = 'baseball sliding shorts men'
$ws.Range('A41').Value = 'knee pads for basketball youth boys'
$ws.Range('A42').Value = 'youth leggings boys basketball'
$ws.Range('A43').Value = 'basketballs leggings'
$ws.Range('A44').Value = 'softball items'
$ws.Range('A45').Value = 'football leggings for men'
$ws.Range('A46').Value = 'basketball compression pants youth'
$ws.Range('A47').Value = 'compression pants men black'
$ws.Range('A48').Value = 'mens work pants with knee pads'
$ws.Range('A49').Value = 'youth tights'
$ws.Range('A50').Value = 'mens down pants'
$ws.Range('A51').Value = 'compression basketball pants youth'
$ws.Range('A52').Value = 'basketball compression pants women'
$ws.Range('A53').Value = 'softball slider'
$ws.Range('A54').Value = 'work knee pads under pants'
$ws.Range('A55').Value = 'calf compression pants'
$ws.Range('A56').Value = 'leg pads basketball'
$ws.Range('A57').Value = 'compression pants men football'
$ws.Range('A58').Value = 'softball shorts men'
$ws.Range('A59').Value = 'mens softball shorts'
$ws.Range('A60').Value = 'mens capri shorts below knee'
$ws.Range('A61').Value = 'compression tights youth'
$ws.Range('A62').Value = 'adidas knee pads'
$ws.Range('A63').Value = 'athletic leggings mens'
$ws.Range('A64').Value = 'knee compression shorts'
$ws.Range('A65').Value = 'compression knee pads pair'
$ws.Range('A66').Value = 'indoor volleyball knee pads'
$ws.Range('A67').Value = 'youth basketball pants boys'
$ws.Range('A68').Value = 'softball sliding shorts girls padded'
$ws.Range('A69').Value = 'lacrosse compression shorts padded'
$ws.Range('A70').Value = 'taken leggings'
$ws.Range('A71').Value = 'mens basketball tights'
$ws.Range('A72').Value = 'black compression pants men'
$ws.Range('A73').Value = 'softball sliding shorts'
$ws.Range('A74').Value = 'girls basketball knee pads youth'
$ws.Range('A75').Value = 'kneepads basketball'
$ws.Range('A76').Value = 'knee protector for construction'
$ws.Range('A77').Value = 'knee pad sleeve basketball'
$ws.Range('A78').Value = 'youth basketball leggings'
$ws.Range('A79').Value = 'medium compression pants'
$ws.Range('A80').Value = 'football knee pads for men'
$ws.Range('A81').Value = 'mens compression leggings'
$ws.Range('A82').Value = 'knee pad for yoga'
$ws.Range('A83').Value = 'boys knee pads basketball'
$ws.Range('A84').Value = 'baseball sliding'
$ws.Range('A85').Value = 'knee pads for biking men'
$ws.Range('A86').Value = 'knee pads girls basketball'
$ws.Range('A87').Value = 'black football pants'
$ws.Range('A88').Value = 'lacrosse pants'
$ws.Range('A89').Value = 'team work softball pants'
$ws.Range('A90').Value = 'long basketball shorts for men below knee'
$ws.Range('A91').Value = 'basketball hex pads'
$ws.Range('A92').Value = 'compression shorts padded basketball'
$ws.Range('A93').Value = 'big boys tights'
$ws.Range('A94').Value = 'basketball tights'
$ws.Range('A95').Value = 'hockey knee pads adult'
$ws.Range('A96').Value = 'padded compression pants football'
$ws.Range('A97').Value = 'hockey hip pads adult'
$ws.Range('A98').Value = 'soccer goalkeeper pads'
$ws.Range('A99').Value = 'athletic mens leggings'
$ws.Range('A100').Value = 'knee pads flexible'
